$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with row 44, which is the "tunnit yht." total row
# (A44 = "tunnit yht.", B44 = SUM(B2:B43)). We need to push that total row
# down by two rows (to row 46) and use the freed-up rows 44 and 45 for two
# new work-log entries, then fix up the SUM range to include them.

# Insert two blank rows right before the current total row; this shifts the
# existing row 44 (with all of its formatting/formula) down to row 46.
$ws.Rows.Item(44).Insert()
$ws.Rows.Item(44).Insert()

# New row 44: copy the date cell formatting used elsewhere in column A
# (centered, m/d/yyyy) onto A44, then fill in the new entry.
$ws.Range("A15").Copy()
$ws.Range("A44").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(44, 1).Value = 44536
$ws.Cells.Item(44, 2).Value = 2
$ws.Cells.Item(44, 3).Value = "ilmettä uusittu, login ja navbar placeholderit lisätty ja komponentit luotu"

# New row 45
$ws.Cells.Item(45, 2).Value = 1
$ws.Cells.Item(45, 3).Value = "css refaktorin, css tiedosto refaktorin todo huomenna"

# Row 46 is now the total row; extend its SUM range to cover the new rows.
$ws.Cells.Item(46, 2).Formula = "=SUM(B2:B45)"
